# Update countries & provincias Spain
# - Refresh a handful of per-country case counters (new day's figures).
# - A few countries that were tied on "Casos totales" got re-ordered by the
#   source feed, which shows up as whole-row swaps for those pairs.
# - Bump the "Datos actualizados" timestamp.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Get-RowVals($r) {
    $vals = @()
    for ($c = 1; $c -le 8; $c++) {
        $vals += ,($ws.Cells.Item($r, $c).Value())
    }
    return $vals
}

function Set-RowVals($r, $vals) {
    for ($c = 1; $c -le 8; $c++) {
        $ws.Cells.Item($r, $c).Value = $vals[$c - 1]
    }
}

function Swap-Rows($r1, $r2) {
    $v1 = Get-RowVals $r1
    $v2 = Get-RowVals $r2
    Set-RowVals $r1 $v2
    Set-RowVals $r2 $v1
}

# --- Refreshed per-country counters ---------------------------------------

# Row 26: Bielorrusia
$ws.Cells.Item(26, 2).Value = 33371
$ws.Cells.Item(26, 3).Value = 945
$ws.Cells.Item(26, 4).Value = 12057
$ws.Cells.Item(26, 5).Value = 21129
$ws.Cells.Item(26, 7).Value = 6
$ws.Cells.Item(26, 8).Value = 185

# Row 99: Eslovenia
$ws.Cells.Item(99, 5).Value = 22
$ws.Cells.Item(99, 7).Value = 1
$ws.Cells.Item(99, 8).Value = 106

# Row 108: Albania
$ws.Cells.Item(108, 2).Value = 969
$ws.Cells.Item(108, 3).Value = 5
$ws.Cells.Item(108, 4).Value = 771
$ws.Cells.Item(108, 5).Value = 167

# --- Tied countries reordered by the source feed ---------------------------

Swap-Rows 197 198   # Santa Lucia <-> Nueva Caledonia
Swap-Rows 209 210   # Montserrat  <-> Seychelles
Swap-Rows 214 215   # Bonaire, San Eustaquio y Saba <-> Sahara Occidental

# --- Timestamp ---------------------------------------------------------

$ws.Range("A1").Value = "Datos actualizados a 21 de Mayo de 2020 a las 12:05"
